$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2744367.2
$ws.Range("I33").Value = 4937156.5
$ws.Range("J33").Value = 3380.6667
$ws.Range("K33").Value = 4937156.5
$ws.Range("L33").Value = 3380.6667
$ws.Range("M33").Value = -4936927.5
$ws.Range("N33").Value = -3838.6667

$ws.Range("H40").Value = 2125.3142
$ws.Range("I40").Value = 1999.8462
$ws.Range("K40").Value = 1999.8462
$ws.Range("M40").Value = -1824.8462

$ws.Range("H76").Value = 3512.652
$ws.Range("I76").Value = 3567.7
$ws.Range("J76").Value = 3470.3076
$ws.Range("K76").Value = 3567.7
$ws.Range("L76").Value = 3470.3076
$ws.Range("M76").Value = -3252.7
$ws.Range("N76").Value = -4100.3076

$ws.Range("H79").Value = 3512.652
$ws.Range("I79").Value = 3567.7
$ws.Range("J79").Value = 3470.3076
$ws.Range("K79").Value = 3567.7
$ws.Range("L79").Value = 3470.3076
$ws.Range("M79").Value = -2475.7
$ws.Range("N79").Value = -5654.3076

$ws.Range("H86").Value = 3787.125
$ws.Range("J86").Value = 3664.6667
$ws.Range("L86").Value = 3664.6667
$ws.Range("N86").Value = -5910.6667

$ws.Range("H89").Value = 3787.125
$ws.Range("J89").Value = 3664.6667
$ws.Range("L89").Value = 18323.3335
$ws.Range("N89").Value = -29555.3335

$ws.Range("H112").Value = 1978.6666
$ws.Range("J112").Value = 2051.4375
$ws.Range("L112").Value = 6154.3125
$ws.Range("N112").Value = -8370.3125

$ws.Range("H115").Value = 251.83333
$ws.Range("I115").Value = 278
$ws.Range("K115").Value = 834
$ws.Range("M115").Value = 733

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6456.3955
$ws.Range("I32").Value = 5795.732
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 5795.732
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5508.732
$ws.Range("N32").Value = -20574

$ws.Range("H61").Value = 17313.363
$ws.Range("J61").Value = 25951.334
$ws.Range("L61").Value = 25951.334
$ws.Range("N61").Value = -26375.334

$ws.Range("H74").Value = 19525.908
$ws.Range("I74").Value = 21292.357
$ws.Range("J74").Value = 16434.625
$ws.Range("K74").Value = 21292.357
$ws.Range("L74").Value = 16434.625
$ws.Range("M74").Value = -20418.357
$ws.Range("N74").Value = -18182.625

$ws.Range("H77").Value = 19525.908
$ws.Range("I77").Value = 21292.357
$ws.Range("J77").Value = 16434.625
$ws.Range("K77").Value = 106461.785
$ws.Range("L77").Value = 82173.125
$ws.Range("M77").Value = -102093.785
$ws.Range("N77").Value = -90909.125

$ws.Range("H130").Value = 40142
$ws.Range("J130").Value = 40142
$ws.Range("L130").Value = 40142
$ws.Range("N130").Value = -50182

$ws.Range("H132").Value = 1528.1364
$ws.Range("I132").Value = 1477.3889
$ws.Range("J132").Value = 1756.5
$ws.Range("K132").Value = 4432.1667
$ws.Range("L132").Value = 5269.5
$ws.Range("M132").Value = -1902.1667
$ws.Range("N132").Value = -10329.5

$ws.Range("H136").Value = 17313.363
$ws.Range("J136").Value = 25951.334
$ws.Range("L136").Value = 77854.00199999999
$ws.Range("N136").Value = -82954.00199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 325.22223
$ws.Range("I22").Value = 314.94116
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 314.94116
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -141.94116
$ws.Range("N22").Value = -846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2954.9722
$ws.Range("I16").Value = 3052.6775
$ws.Range("K16").Value = 3052.6775
$ws.Range("M16").Value = -2765.6775

$ws.Range("H62").Value = 90231.164
$ws.Range("I62").Value = 203161.2
$ws.Range("J62").Value = 9566.857
$ws.Range("K62").Value = 203161.2
$ws.Range("L62").Value = 9566.857
$ws.Range("M62").Value = -202537.2
$ws.Range("N62").Value = -10814.857

$ws.Range("H65").Value = 90231.164
$ws.Range("I65").Value = 203161.2
$ws.Range("J65").Value = 9566.857
$ws.Range("K65").Value = 1015806
$ws.Range("L65").Value = 47834.285
$ws.Range("M65").Value = -1012686
$ws.Range("N65").Value = -54074.285

$ws.Range("H113").Value = 2954.9722
$ws.Range("I113").Value = 3052.6775
$ws.Range("K113").Value = 3052.6775
$ws.Range("M113").Value = -882.6774999999998

$ws.Range("H132").Value = 26722.584
$ws.Range("I132").Value = 15692.733
$ws.Range("K132").Value = 47078.199
$ws.Range("M132").Value = -44548.199

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 56.24
$ws.Range("I2").Value = 46.387096
$ws.Range("K2").Value = 278.322576
$ws.Range("M2").Value = -165.322576

$ws.Range("H7").Value = 41666732
$ws.Range("I7").Value = 62500064
$ws.Range("J7").Value = 64.75
$ws.Range("K7").Value = 187500192
$ws.Range("L7").Value = 194.25
$ws.Range("M7").Value = -187500080
$ws.Range("N7").Value = -418.25

$ws.Range("H55").Value = 200005020
$ws.Range("J55").Value = 8316.333000000001
$ws.Range("L55").Value = 24948.999
$ws.Range("N55").Value = -25302.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 49995
$ws.Range("J52").Value = 49995
$ws.Range("L52").Value = 49995
$ws.Range("N52").Value = -50513

$ws.Range("H97").Value = 872.55554
$ws.Range("I97").Value = 655.2
$ws.Range("J97").Value = 1144.25
$ws.Range("K97").Value = 655.2
$ws.Range("L97").Value = 1144.25
$ws.Range("M97").Value = -159.2
$ws.Range("N97").Value = -2136.25

$ws.Range("H132").Value = 19473.691
$ws.Range("I132").Value = 32380.143
$ws.Range("K132").Value = 97140.429
$ws.Range("M132").Value = -94610.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1687.25
$ws.Range("I46").Value = 928.5
$ws.Range("K46").Value = 928.5
$ws.Range("M46").Value = -740.5

$ws.Range("H68").Value = 3336.818
$ws.Range("I68").Value = 2861.2727
$ws.Range("J68").Value = 3812.3635
$ws.Range("K68").Value = 2861.2727
$ws.Range("L68").Value = 3812.3635
$ws.Range("M68").Value = -2112.2727
$ws.Range("N68").Value = -5310.363499999999

$ws.Range("H71").Value = 3336.818
$ws.Range("I71").Value = 2861.2727
$ws.Range("J71").Value = 3812.3635
$ws.Range("K71").Value = 14306.3635
$ws.Range("L71").Value = 19061.8175
$ws.Range("M71").Value = -10562.3635
$ws.Range("N71").Value = -26549.8175

$ws.Range("H93").Value = 4817.077
$ws.Range("I93").Value = 5614.8125
$ws.Range("J93").Value = 3540.7
$ws.Range("K93").Value = 5614.8125
$ws.Range("L93").Value = 3540.7
$ws.Range("M93").Value = -4366.8125
$ws.Range("N93").Value = -6036.7

$ws.Range("H99").Value = 33181.5
$ws.Range("I99").Value = 31312.777
$ws.Range("K99").Value = 31312.777
$ws.Range("M99").Value = -28317.777

$ws.Range("H128").Value = 89142.664
$ws.Range("J128").Value = 89142.664
$ws.Range("L128").Value = 89142.664
$ws.Range("N128").Value = -99102.664

$ws.Range("H130").Value = 90000
$ws.Range("J130").Value = 90000
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040

$ws.Range("H131").Value = 767500
$ws.Range("J131").Value = 767500
$ws.Range("L131").Value = 767500
$ws.Range("N131").Value = -777580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 34983.668
$ws.Range("I61").Value = 34983.668
$ws.Range("K61").Value = 34983.668
$ws.Range("M61").Value = -34691.668

$ws.Range("H132").Value = 142132.83
$ws.Range("I132").Value = 200866.28
$ws.Range("J132").Value = 29016.555
$ws.Range("K132").Value = 602598.84
$ws.Range("L132").Value = 87049.66500000001
$ws.Range("M132").Value = -600068.84
$ws.Range("N132").Value = -92109.66500000001
